$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before column H. This shifts the old
#    H/I/J (Notes/Type/Rule For) columns one position to the right,
#    including their data validations.
$ws.Columns("H").Insert()

# 2. Update the header text for the (now widened) Amount column and the
#    brand-new Folio-currency amount column.
$ws.Range("G1").Value = "Amount (Fund Currency)*"
$ws.Range("H1").Value = "Amount (Folio Currency)"

# 3. The cell comments did not travel with the column insert (they are
#    still anchored to H1/I1/J1, which now hold the Notes/Type/Rule For
#    comments that logically belong one column to the right at I1/J1/K1).
#    Capture the existing comment text, then delete and re-add each one
#    at its correct new home.
$notesComment = $ws.Range("H1").Comment.Text()
$typeComment = $ws.Range("I1").Comment.Text()
$ruleForComment = $ws.Range("J1").Comment.Text()

$ws.Range("H1").Comment.Delete()
$ws.Range("I1").Comment.Delete()
$ws.Range("J1").Comment.Delete()

$ws.Range("I1").AddComment($notesComment)
$ws.Range("J1").AddComment($typeComment)
$ws.Range("K1").AddComment($ruleForComment)

# 4. Match the final active selection recorded in the saved workbook.
$ws.Range("H1").Select()

"done"
